# ChatGPT_Results.xlsx update:
#  - A4/A5 "(x = 0) & (y = 0)" -> "(x = 10) & (y = 10)"
#  - Row 6 and Row 8 formatting de-duplicated (drop the redundant
#    applyFill variants of the left/center and center/center wrap
#    styles, and the redundant default style) so those rows use the
#    same cell styles as the rest of the table.
#  - Row 4 height grows to an explicit 88.2pt custom height.
#  - The active view scrolls back to the top and the selection moves
#    to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expression text update (rows 4 & 5, column A) ---------------------
$ws.Range("A4").Value = "(x = 10) & (y = 10)"
$ws.Range("A5").Value = "(x = 10) & (y = 10)"

# --- De-duplicate the one-off "applyFill" styles used by rows 6 & 8 ----
# These xf records are visually identical to the plain (non-fill) styles
# already used everywhere else (fillId is "none" in both cases), so
# re-applying the equivalent alignment collapses the cells back onto the
# shared style entries instead of the redundant ones.
foreach ($addr in @("A6", "B6", "D6", "A8", "B8", "D8")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4131   # xlLeft
    $c.VerticalAlignment = -4108     # xlCenter
    $c.WrapText = $true
}

foreach ($addr in @("C6", "E6", "F6", "G6", "C8", "E8", "F8", "G8")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4108     # xlCenter
    $c.WrapText = $true
}

$ws.Range("H6").Style = "Normal"
$ws.Range("H8").Style = "Normal"

# --- Row 4 explicit custom height --------------------------------------
$ws.Rows("4").RowHeight = 88.2

# --- View state: scroll back to top, select A5 -------------------------
# (topLeftCell/pane-scroll state isn't exposed by this host's object
# model, so only the selection itself is settable here.)
$ws.Range("A5").Select() | Out-Null
